$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-02-06 Tuesday"; New = "2024-02-07 Wednesday" },
    @{ Old = "164×7=1148"; New = "276×7=1932" },
    @{ Old = "718×5=3590"; New = "856×2=1712" },
    @{ Old = "969×7=6783"; New = "859×5=4295" },
    @{ Old = "927×4=3708"; New = "437×8=3496" },
    @{ Old = "493×4=1972"; New = "589×9=5301" },
    @{ Old = "579×9=5211"; New = "181×4=724" },
    @{ Old = "936×2=1872"; New = "335×3=1005" },
    @{ Old = "518×5=2590"; New = "485×5=2425" },
    @{ Old = "929×9=8361"; New = "516×8=4128" },
    @{ Old = "303×6=1818"; New = "473×5=2365" },
    @{ Old = "757×4=3028"; New = "798×4=3192" },
    @{ Old = "384×3=1152"; New = "191×7=1337" },
    @{ Old = "210×6=1260"; New = "739×6=4434" },
    @{ Old = "971×5=4855"; New = "882×9=7938" },
    @{ Old = "220×8=1760"; New = "233×3=699" },
    @{ Old = "808×2=1616"; New = "792×9=7128" },
    @{ Old = "747×7=5229"; New = "297×2=594" },
    @{ Old = "588×6=3528"; New = "858×5=4290" },
    @{ Old = "470×3=1410"; New = "564×3=1692" },
    @{ Old = "651×8=5208"; New = "246×3=738" },
    @{ Old = "776×2=1552"; New = "908×5=4540" },
    @{ Old = "305×2=610"; New = "320×9=2880" },
    @{ Old = "262×8=2096"; New = "285×6=1710" },
    @{ Old = "983×6=5898"; New = "984×4=3936" },
    @{ Old = "189×3=567"; New = "882×2=1764" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}

$d.Save()
